# Auto-generated edit script reproducing the diff against Sheets/Ultros_Profits.xlsx
# Updates computed price/profit columns (H:N) for specific rows across all 8 item-category sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12:N12").ClearContents()
$ws.Range("H40").Value = 7164.3335
$ws.Range("I40").Value = 1993.5
$ws.Range("J40").Value = 9749.75
$ws.Range("K40").Value = 1993.5
$ws.Range("L40").Value = 9749.75
$ws.Range("M40").Value = -1818.5
$ws.Range("N40").Value = -10099.75
$ws.Range("H43").Value = 4510.4707
$ws.Range("J43").Value = 4045.2
$ws.Range("L43").Value = 4045.2
$ws.Range("N43").Value = -4183.2
$ws.Range("H127").Value = 4892.6665
$ws.Range("I127").Value = 1610.1818
$ws.Range("K127").Value = 4830.5454
$ws.Range("M127").Value = 129.4546

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 2019.7
$ws.Range("I19").Value = 889.75
$ws.Range("K19").Value = 889.75
$ws.Range("M19").Value = -660.75
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29:N29").ClearContents()
$ws.Range("H45").Value = 6123.8335
$ws.Range("I45").Value = 5442
$ws.Range("J45").Value = 7487.5
$ws.Range("K45").Value = 5442
$ws.Range("L45").Value = 7487.5
$ws.Range("M45").Value = -5065
$ws.Range("N45").Value = -8241.5
$ws.Range("H56").Value = 3200
$ws.Range("J56").Value = 8000
$ws.Range("L56").Value = 8000
$ws.Range("N56").Value = -9484
$ws.Range("H104").Value = 16725
$ws.Range("J104").Value = 16725
$ws.Range("L104").Value = 16725
$ws.Range("N104").Value = -23713
$ws.Range("H110").Value = 5795.846
$ws.Range("I110").Value = 5693.3
$ws.Range("K110").Value = 5693.3
$ws.Range("M110").Value = -3648.3
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 9407.862999999999
$ws.Range("I107").Value = 9007.941000000001
$ws.Range("K107").Value = 9007.941000000001
$ws.Range("M107").Value = -7087.941000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 355.75
$ws.Range("I22").Value = 440.4
$ws.Range("K22").Value = 440.4
$ws.Range("M22").Value = -90.39999999999998
$ws.Range("H31").Value = 3021.7273
$ws.Range("I31").Value = 1850.5
$ws.Range("K31").Value = 1850.5
$ws.Range("M31").Value = -1555.5
$ws.Range("H34").Value = 3021.7273
$ws.Range("I34").Value = 1850.5
$ws.Range("K34").Value = 1850.5
$ws.Range("M34").Value = -1648.5
$ws.Range("H86").Value = 35306.535
$ws.Range("I86").Value = 45462
$ws.Range("K86").Value = 45462
$ws.Range("M86").Value = -44339
$ws.Range("H89").Value = 35306.535
$ws.Range("I89").Value = 45462
$ws.Range("K89").Value = 227310
$ws.Range("M89").Value = -221694
$ws.Range("H132").Value = 3437.625
$ws.Range("I132").Value = 2897.4
$ws.Range("J132").Value = 4338
$ws.Range("K132").Value = 8692.200000000001
$ws.Range("L132").Value = 13014
$ws.Range("M132").Value = -6162.200000000001
$ws.Range("N132").Value = -18074

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1385.7
$ws.Range("I98").Value = 1596.8182
$ws.Range("J98").Value = 1127.6666
$ws.Range("K98").Value = 4790.4546
$ws.Range("L98").Value = 3382.9998
$ws.Range("M98").Value = -3292.4546
$ws.Range("N98").Value = -6378.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 159.08333
$ws.Range("I2").Value = 55.384617
$ws.Range("J2").Value = 281.63635
$ws.Range("K2").Value = 55.384617
$ws.Range("L2").Value = 281.63635
$ws.Range("M2").Value = 57.615383
$ws.Range("N2").Value = -507.63635
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H43").Value = 27391.217
$ws.Range("I43").Value = 17499.8
$ws.Range("J43").Value = 35000
$ws.Range("K43").Value = 17499.8
$ws.Range("L43").Value = 35000
$ws.Range("M43").Value = -17348.8
$ws.Range("N43").Value = -35302
$ws.Range("H70").Value = 146169.12
$ws.Range("I70").Value = 285876.5
$ws.Range("K70").Value = 285876.5
$ws.Range("M70").Value = -285606.5
$ws.Range("H73").Value = 146169.12
$ws.Range("I73").Value = 285876.5
$ws.Range("K73").Value = 285876.5
$ws.Range("M73").Value = -284940.5
$ws.Range("H99").Value = 25951
$ws.Range("J99").Value = 44495
$ws.Range("L99").Value = 44495
$ws.Range("N99").Value = -48987
$ws.Range("H102").Value = 7964.8887
$ws.Range("I102").Value = 6857.143
$ws.Range("K102").Value = 6857.143
$ws.Range("M102").Value = -5235.143
$ws.Range("H113").Value = 6187.9
$ws.Range("I113").Value = 2868.4285
$ws.Range("K113").Value = 2868.4285
$ws.Range("M113").Value = -698.4285
$ws.Range("H122").Value = 4848.5
$ws.Range("I122").Value = 4103.6113
$ws.Range("K122").Value = 12310.8339
$ws.Range("M122").Value = -9860.833899999998
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126:N126").ClearContents()
$ws.Range("H132").Value = 9981.291999999999
$ws.Range("I132").Value = 9264.764999999999
$ws.Range("J132").Value = 11721.429
$ws.Range("K132").Value = 27794.295
$ws.Range("L132").Value = 35164.287
$ws.Range("M132").Value = -25264.295
$ws.Range("N132").Value = -40224.287

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 22727896
$ws.Range("J22").Value = 1993
$ws.Range("K22").Value = 22727896
$ws.Range("L22").Value = 1993
$ws.Range("M22").Value = -22727601
$ws.Range("N22").Value = -2583
$ws.Range("H25").Value = 199999
$ws.Range("I25").Value = 199999
$ws.Range("K25").Value = 199999
$ws.Range("M25").Value = -199769
$ws.Range("I27").Value = 22727896
$ws.Range("J27").Value = 1993
$ws.Range("K27").Value = 22727896
$ws.Range("L27").Value = 1993
$ws.Range("M27").Value = -22727789
$ws.Range("N27").Value = -2207
$ws.Range("H36").Value = 87499.5
$ws.Range("J36").Value = 87499.5
$ws.Range("L36").Value = 87499.5
$ws.Range("N36").Value = -88623.5
$ws.Range("H40").Value = 68099.75
$ws.Range("I40").Value = 24133
$ws.Range("K40").Value = 24133
$ws.Range("M40").Value = -23997
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H94").Value = 52499.75
$ws.Range("J94").Value = 52499.75
$ws.Range("L94").Value = 52499.75
$ws.Range("N94").Value = -53851.75
$ws.Range("H99").Value = 200285
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 200285
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 200285
$ws.Range("N99").Value = -206275
$ws.Range("M99").ClearContents()
$ws.Range("H100").Value = 112460
$ws.Range("I100").Value = 140388.75
$ws.Range("K100").Value = 140388.75
$ws.Range("M100").Value = -139847.75
$ws.Range("H116").Value = 246793.8
$ws.Range("J116").Value = 246793.8
$ws.Range("L116").Value = 246793.8
$ws.Range("N116").Value = -255971.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 3999.5
$ws.Range("I8").Value = 4000
$ws.Range("J8").Value = 3999.3333
$ws.Range("K8").Value = 4000
$ws.Range("L8").Value = 3999.3333
$ws.Range("M8").Value = -3860
$ws.Range("N8").Value = -4279.3333
$ws.Range("H122").Value = 1659.875
$ws.Range("I122").Value = 1611.2858
$ws.Range("K122").Value = 4833.857400000001
$ws.Range("M122").Value = -2383.857400000001
